# Update odds values in row 3 (match: Santa Fe vs Chico) to reflect the
# latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value  = 1.55
$ws.Range("H3").Value  = 3.8
$ws.Range("I3").Value  = 6.25
$ws.Range("N3").Value  = 7.5
$ws.Range("W3").Value  = 5.5
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 7.5
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 81
$ws.Range("AK3").Value = 67
$ws.Range("AM3").Value = 51
$ws.Range("AW3").Value = 7.5
